{"js": "const pairs = [\n  [\"2024-09-19 Thursday\", \"2024-09-20 Friday\"],\n  [\"58+7=65\", \"16+2=18\"],\n  [\"5+62=67\", \"44+24=68\"],\n  [\"54+43=97\", \"11+50=61\"],\n  [\"12+73=85\", \"48+23=71\"],\n  [\"40+27=67\", \"96+0=96\"],\n  [\"13+39=52\", \"74-74=0\"],\n  [\"61+7=68\", \"34+59=93\"],\n  [\"10+5=15\", \"66-24=42\"],\n  [\"66+27=93\", \"64-61=3\"],\n  [\"49-6=43\", \"7-3=4\"],\n  [\"57-43=14\", \"35-21=14\"],\n  [\"43+37=80\", \"8+3=11\"],\n  [\"18+44=62\", \"89-17=72\"],\n  [\"33+39=72\", \"0+61=61\"],\n  [\"38+1=39\", \"56+2=58\"],\n  [\"55+34=89\", \"1+23=24\"],\n  [\"67-0=67\", \"61-29=32\"],\n  [\"66+22=88\", \"31+57=88\"],\n  [\"13+20=33\", \"72-43=29\"],\n  [\"21+5=26\", \"23+57=80\"],\n  [\"57-23=34\", \"14+3=17\"],\n  [\"96-25=71\", \"78-50=28\"],\n  [\"60-48=12\", \"86-49=37\"],\n  [\"27+28=55\", \"52-49=3\"],\n  [\"7+0=7\", \"74+4=78\"],\n  [\"99-33=66\", \"90-11=79\"],\n  [\"54-12=42\", \"1+78=79\"],\n  [\"7+31=38\", \"6+88=94\"],\n  [\"78+20=98\", \"74-12=62\"],\n  [\"41-29=12\", \"20+35=55\"],\n  [\"53-1=52\", \"9+60=69\"],\n  [\"95-8=87\", \"88-46=42\"],\n  [\"39-10=29\", \"11+50=61\"],\n  [\"64-21=43\", \"73-16=57\"],\n  [\"0+26=26\", \"95-46=49\"],\n  [\"57-11=46\", \"76-4=72\"],\n  [\"53+33=86\", \"15+25=40\"],\n  [\"9+6=15\", \"19+27=46\"],\n  [\"97-58=39\", \"29+11=40\"],\n  [\"98-76=22\", \"12+31=43\"],\n  [\"63-30=33\", \"28+26=54\"],\n  [\"92-0=92\", \"50+24=74\"],\n  [\"26-12=14\", \"80-40=40\"],\n  [\"36-12=24\", \"28+32=60\"],\n  [\"82-79=3\", \"5+54=59\"],\n  [\"80-30=50\", \"90-77=13\"],\n  [\"36+43=79\", \"25-3=22\"],\n  [\"39+48=87\", \"32+36=68\"],\n  [\"92-61=31\", \"10-3=7\"],\n  [\"24+61=85\", \"47-40=7\"],\n  [\"94-61=33\", \"87+4=91\"],\n  [\"33+62=95\", \"21+42=63\"],\n  [\"5+45=50\", \"9+0=9\"],\n  [\"73-25=48\", \"47+2=49\"],\n  [\"84-74=10\", \"13+7=20\"],\n  [\"46+24=70\", \"97-31=66\"],\n  [\"52-17=35\", \"99-52=47\"],\n  [\"50-31=19\", \"33+2=35\"],\n  [\"8+80=88\", \"61-33=28\"],\n  [\"26+59=85\", \"20+24=44\"],\n  [\"1+29=30\", \"84+15=99\"],\n  [\"8+49=57\", \"31-23=8\"],\n  [\"58-12=46\", \"35+27=62\"],\n  [\"40+59=99\", \"91-91=0\"],\n  [\"11+20=31\", \"44+36=80\"],\n  [\"3+60=63\", \"18+28=46\"],\n  [\"46+43=89\", \"95-17=78\"],\n  [\"74-31=43\", \"35-12=23\"],\n  [\"46+27=73\", \"57-45=12\"],\n  [\"83+10=93\", \"72-26=46\"],\n  [\"30-6=24\", \"5+19=24\"],\n  [\"76+13=89\", \"4+68=72\"],\n  [\"82-65=17\", \"57+13=70\"],\n  [\"17-2=15\", \"78+11=89\"],\n  [\"31+58=89\", \"2+7=9\"],\n  [\"69-38=31\", \"36+8=44\"],\n  [\"38+60=98\", \"20+44=64\"],\n  [\"91-2=89\", \"27+13=40\"],\n  [\"79+2=81\", \"39-23=16\"],\n  [\"52-20=32\", \"51-46=5\"],\n  [\"81-68=13\", \"25+53=78\"],\n  [\"91+5=96\", \"13+13=26\"],\n  [\"68-55=13\", \"70-49=21\"],\n  [\"40+26=66\", \"44-3=41\"],\n  [\"41+10=51\", \"42-33=9\"],\n  [\"33-0=33\", \"95-53=42\"],\n  [\"17-1=16\", \"23+17=40\"],\n  [\"7+23=30\", \"54-48=6\"],\n  [\"0+86=86\", \"56-56=0\"],\n  [\"61+4=65\", \"90+3=93\"],\n  [\"0+57=57\", \"96-61=35\"],\n  [\"83-55=28\", \"88-18=70\"],\n  [\"47-22=25\", \"44+1=45\"],\n  [\"29+0=29\", \"17+11=28\"],\n  [\"13+23=36\", \"2+87=89\"],\n  [\"56-23=33\", \"0+9=9\"],\n  [\"97-2=95\", \"25+63=88\"],\n  [\"89-57=32\", \"59+30=89\"],\n  [\"15+65=80\", \"99-34=65\"],\n  [\"18+59=77\", \"29+52=81\"],\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Not found: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  totalReplaced++;\n}\nawait context.sync();\nreturn 'replaced: ' + totalReplaced;", "ps1": "$pairs = @(\n    @(\"2024-09-19 Thursday\", \"2024-09-20 Friday\"),\n    @(\"58+7=65\", \"16+2=18\"),\n    @(\"5+62=67\", \"44+24=68\"),\n    @(\"54+43=97\", \"11+50=61\"),\n    @(\"12+73=85\", \"48+23=71\"),\n    @(\"40+27=67\", \"96+0=96\"),\n    @(\"13+39=52\", \"74-74=0\"),\n    @(\"61+7=68\", \"34+59=93\"),\n    @(\"10+5=15\", \"66-24=42\"),\n    @(\"66+27=93\", \"64-61=3\"),\n    @(\"49-6=43\", \"7-3=4\"),\n    @(\"57-43=14\", \"35-21=14\"),\n    @(\"43+37=80\", \"8+3=11\"),\n    @(\"18+44=62\", \"89-17=72\"),\n    @(\"33+39=72\", \"0+61=61\"),\n    @(\"38+1=39\", \"56+2=58\"),\n    @(\"55+34=89\", \"1+23=24\"),\n    @(\"67-0=67\", \"61-29=32\"),\n    @(\"66+22=88\", \"31+57=88\"),\n    @(\"13+20=33\", \"72-43=29\"),\n    @(\"21+5=26\", \"23+57=80\"),\n    @(\"57-23=34\", \"14+3=17\"),\n    @(\"96-25=71\", \"78-50=28\"),\n    @(\"60-48=12\", \"86-49=37\"),\n    @(\"27+28=55\", \"52-49=3\"),\n    @(\"7+0=7\", \"74+4=78\"),\n    @(\"99-33=66\", \"90-11=79\"),\n    @(\"54-12=42\", \"1+78=79\"),\n    @(\"7+31=38\", \"6+88=94\"),\n    @(\"78+20=98\", \"74-12=62\"),\n    @(\"41-29=12\", \"20+35=55\"),\n    @(\"53-1=52\", \"9+60=69\"),\n    @(\"95-8=87\", \"88-46=42\"),\n    @(\"39-10=29\", \"11+50=61\"),\n    @(\"64-21=43\", \"73-16=57\"),\n    @(\"0+26=26\", \"95-46=49\"),\n    @(\"57-11=46\", \"76-4=72\"),\n    @(\"53+33=86\", \"15+25=40\"),\n    @(\"9+6=15\", \"19+27=46\"),\n    @(\"97-58=39\", \"29+11=40\"),\n    @(\"98-76=22\", \"12+31=43\"),\n    @(\"63-30=33\", \"28+26=54\"),\n    @(\"92-0=92\", \"50+24=74\"),\n    @(\"26-12=14\", \"80-40=40\"),\n    @(\"36-12=24\", \"28+32=60\"),\n    @(\"82-79=3\", \"5+54=59\"),\n    @(\"80-30=50\", \"90-77=13\"),\n    @(\"36+43=79\", \"25-3=22\"),\n    @(\"39+48=87\", \"32+36=68\"),\n    @(\"92-61=31\", \"10-3=7\"),\n    @(\"24+61=85\", \"47-40=7\"),\n    @(\"94-61=33\", \"87+4=91\"),\n    @(\"33+62=95\", \"21+42=63\"),\n    @(\"5+45=50\", \"9+0=9\"),\n    @(\"73-25=48\", \"47+2=49\"),\n    @(\"84-74=10\", \"13+7=20\"),\n    @(\"46+24=70\", \"97-31=66\"),\n    @(\"52-17=35\", \"99-52=47\"),\n    @(\"50-31=19\", \"33+2=35\"),\n    @(\"8+80=88\", \"61-33=28\"),\n    @(\"26+59=85\", \"20+24=44\"),\n    @(\"1+29=30\", \"84+15=99\"),\n    @(\"8+49=57\", \"31-23=8\"),\n    @(\"58-12=46\", \"35+27=62\"),\n    @(\"40+59=99\", \"91-91=0\"),\n    @(\"11+20=31\", \"44+36=80\"),\n    @(\"3+60=63\", \"18+28=46\"),\n    @(\"46+43=89\", \"95-17=78\"),\n    @(\"74-31=43\", \"35-12=23\"),\n    @(\"46+27=73\", \"57-45=12\"),\n    @(\"83+10=93\", \"72-26=46\"),\n    @(\"30-6=24\", \"5+19=24\"),\n    @(\"76+13=89\", \"4+68=72\"),\n    @(\"82-65=17\", \"57+13=70\"),\n    @(\"17-2=15\", \"78+11=89\"),\n    @(\"31+58=89\", \"2+7=9\"),\n    @(\"69-38=31\", \"36+8=44\"),\n    @(\"38+60=98\", \"20+44=64\"),\n    @(\"91-2=89\", \"27+13=40\"),\n    @(\"79+2=81\", \"39-23=16\"),\n    @(\"52-20=32\", \"51-46=5\"),\n    @(\"81-68=13\", \"25+53=78\"),\n    @(\"91+5=96\", \"13+13=26\"),\n    @(\"68-55=13\", \"70-49=21\"),\n    @(\"40+26=66\", \"44-3=41\"),\n    @(\"41+10=51\", \"42-33=9\"),\n    @(\"33-0=33\", \"95-53=42\"),\n    @(\"17-1=16\", \"23+17=40\"),\n    @(\"7+23=30\", \"54-48=6\"),\n    @(\"0+86=86\", \"56-56=0\"),\n    @(\"61+4=65\", \"90+3=93\"),\n    @(\"0+57=57\", \"96-61=35\"),\n    @(\"83-55=28\", \"88-18=70\"),\n    @(\"47-22=25\", \"44+1=45\"),\n    @(\"29+0=29\", \"17+11=28\"),\n    @(\"13+23=36\", \"2+87=89\"),\n    @(\"56-23=33\", \"0+9=9\"),\n    @(\"97-2=95\", \"25+63=88\"),\n    @(\"89-57=32\", \"59+30=89\"),\n    @(\"15+65=80\", \"99-34=65\"),\n    @(\"18+59=77\", \"29+52=81\"),\n)\n\n$d = $word.ActiveDocument\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $result = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $result) {\n        throw \"Replace failed for: $oldText\"\n    }\n}"}
